# feat: parse `Edm.Time` when excel data is text (#573)
#
# The "time" column (H) held numeric Excel time serial values. With this
# change, a time coming from the OData/Edm.Time source is written to the
# sheet as plain text (e.g. "16:00:00") instead of a numeric time value, so
# the cell needs a Text ("@") number format and a literal string value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("time[time]" column H) now carries the Edm.Time value as text
# rather than a numeric Excel time serial.
$cell = $ws.Range("H3")
$cell.NumberFormat = "@"
$cell.Value = "16:00:00"

# Reflect the author's last selection in the saved sheet view.
$ws.Range("H13").Select()
